# Apply updated crypto price/volume figures to sheet1 (rows 2-51, columns D & E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.977.49"
$ws.Range("E2").Value = "  -1.69%  "
$ws.Range("D3").Value = "1.820.52"
$ws.Range("E3").Value = "  -1.22%  "
$ws.Range("E4").Value = "  -0.57%  "
$ws.Range("D5").Value = "'309.82"
$ws.Range("E5").Value = "  -1.87%  "
$ws.Range("D6").Value = "'1.009"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").Value = "'0.4631"
$ws.Range("E7").Value = "  -2.97%  "
$ws.Range("D8").Value = "'0.3639"
$ws.Range("E8").Value = "  -1.91%  "
$ws.Range("D9").Value = "'0.07291"
$ws.Range("E9").Value = "  -2.45%  "
$ws.Range("D10").Value = "'0.8662"
$ws.Range("E10").Value = "  -2.44%  "
$ws.Range("D11").Value = "'19.81"
$ws.Range("E11").Value = "  -3.38%  "
$ws.Range("D12").Value = "1.882.78"
$ws.Range("E12").Value = "  +1.97%  "
$ws.Range("D13").Value = "'0.07601"
$ws.Range("E13").Value = "  +3.02%  "
$ws.Range("D14").Value = "'93.16"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").Value = "'5.328"
$ws.Range("E15").Value = "  -2.98%  "
$ws.Range("D16").Value = "'6.481"
$ws.Range("E16").Value = "  -1.89%  "
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").Value = "'0.000008630"
$ws.Range("E18").Value = "  -2.56%  "
$ws.Range("D19").Value = "'1.008"
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").Value = "27.363.80"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("E21").Value = "  -2.52%  "
$ws.Range("D22").Value = "'5.166"
$ws.Range("E22").Value = "  -3.62%  "
$ws.Range("E23").Value = "  -1.61%  "
$ws.Range("D24").Value = "2.116.26"
$ws.Range("E24").Value = "  +1.81%  "
$ws.Range("D25").Value = "'151.88"
$ws.Range("E25").Value = "  -0.60%  "
$ws.Range("D26").Value = "'1.854"
$ws.Range("E26").Value = "  -2.68%  "
$ws.Range("E27").Value = "  -2.30%  "
$ws.Range("D28").Value = "'2.096"
$ws.Range("E28").Value = "  -3.53%  "
$ws.Range("D29").Value = "'5.086"
$ws.Range("E29").Value = "  -3.72%  "
$ws.Range("D30").Value = "'115.74"
$ws.Range("E30").Value = "  -2.11%  "
$ws.Range("D31").Value = "'0.08907"
$ws.Range("E31").Value = "  -0.94%  "
$ws.Range("D32").Value = "'2.951"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").Value = "'0.7288"
$ws.Range("E33").Value = "  -4.22%  "
$ws.Range("D34").Value = "'1.139"
$ws.Range("E34").Value = "  -3.36%  "
$ws.Range("D35").Value = "'4.430"
$ws.Range("E35").Value = "  -3.07%  "
$ws.Range("E36").Value = "  -0.57%  "
$ws.Range("D37").Value = "'2.509"
$ws.Range("E37").Value = "  +5.10%  "
$ws.Range("D38").Value = "'0.05277"
$ws.Range("E38").Value = "  -1.84%  "
$ws.Range("E39").Value = "  -2.97%  "
$ws.Range("D40").Value = "'0.01916"
$ws.Range("E40").Value = "  -2.82%  "
$ws.Range("D41").Value = "'2.932"
$ws.Range("E41").Value = "  -2.33%  "
$ws.Range("D42").Value = "'7.167"
$ws.Range("E42").Value = "  -2.19%  "
$ws.Range("D43").Value = "'0.5223"
$ws.Range("E43").Value = "  -2.76%  "
$ws.Range("D44").Value = "'0.1633"
$ws.Range("E44").Value = "  -2.22%  "
$ws.Range("D45").Value = "'8.264"
$ws.Range("E45").Value = "  -3.61%  "
$ws.Range("D46").Value = "'0.4867"
$ws.Range("E46").Value = "  -2.49%  "
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("D48").Value = "'10.11"
$ws.Range("E48").Value = "  -4.62%  "
$ws.Range("D49").Value = "'103.25"
$ws.Range("E49").Value = "  -1.94%  "
$ws.Range("D50").Value = "'1.634"
$ws.Range("E50").Value = "  -3.12%  "
$ws.Range("D51").Value = "'0.06224"
$ws.Range("E51").Value = "  -1.65%  "
